# Add the 2020 column (O) of data, extending the existing 2009-2019 (D:N) series.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last existing data column (N) onto the new column (O)
# for both the header row (4) and the data row (5), then fill in the new values.
$ws.Range("N4:N5").Copy() | Out-Null
$ws.Range("O4:O5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

$ws.Range("O4").Value = 2020
$ws.Range("O5").Value = 83.3

$excel.CutCopyMode = 0

# Move the active selection, matching the saved view state of the edited file.
$ws.Range("O12").Select() | Out-Null
